$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price (D) and 1h volume-change (E) columns with the
# latest feed snapshot. Price cells whose new reading parses as a plain
# number are forced back to text (matches the columns existing text
# storage, e.g. "1.00"/"587.13") by setting NumberFormat to "@" first.

$ws.Range("D2").Value = "64.797.35"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.509.68"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.13"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.35"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "3.508.84"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "4.108.62"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.77"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").Value = "3.515.10"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "64.817.80"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.27"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.33"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.577"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "3.653.14"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.25"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").Value = "3.516.29"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.96"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.75"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.20"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.98"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0809"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.41"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.30"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "2.486.16"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.88"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.902"
$ws.Range("E51").Value = "  +3.58%  "
